# Slide 7 ("LOSS FUNCTION" slide), shape "TextBox 8": the bullet describing
# Adam is edited from
#   "Adam combines SGD and RMSProp, which use momentum and adaptive learning rates"
# to
#   "Adam combines Momentum and RMSProp, which use momentum terms and adaptive learning rates"
#
# i.e. "SGD" -> "Momentum" and "momentum" -> "momentum terms". PowerPoint splits
# the paragraph's single run into several runs at the edited word boundaries
# (each edited/typed word keeps its own run), which we reproduce here by
# re-assigning .Text on successive Characters() sub-ranges so the run
# formatting (color/typeface/size) is inherited/copied onto every new run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item("TextBox 8")
$tr = $shp.TextFrame.TextRange

$tr.Characters(1, 14).Text  = "Adam combines "
$tr.Characters(15, 4).Text  = "Momentum "
$tr.Characters(24, 4).Text  = "and "
$tr.Characters(28, 7).Text  = "RMSProp"
$tr.Characters(35, 12).Text = ", which use "
$tr.Characters(47, 9).Text  = "momentum terms "
$tr.Characters(62, 27).Text = "and adaptive learning rates"

Write-Output $tr.Text
